# Calculator Keyboard Layout update
# Commit: "Add processing for absolute value" - repurpose the unused
# comparison-operator keys (<, >, <=, >=) on the secondary keypad (columns
# L/M and W/X of row 4) into an absolute-value key ("|x|") plus matching
# parenthesis keys, mirroring the primary keypad's "(" / ")" keys.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# L4 was "<" -> becomes "|x|", restyled like the other special-function
# keys on that keypad (e.g. R4 "Pi": size-18 Consolas on yellow fill).
$ws.Range("R4").Copy() | Out-Null
$ws.Range("L4").PasteSpecial(-4122) | Out-Null
$ws.Range("L4").Value = "|x|"

# M4 was ">" -> becomes ")", restyled to match the primary keypad's
# parenthesis keys (A4/B4 style).
$ws.Range("A4").Copy() | Out-Null
$ws.Range("M4").PasteSpecial(-4122) | Out-Null
$ws.Range("M4").Value = ")"

# W4 was "<=" -> becomes "(", same parenthesis style.
$ws.Range("A4").Copy() | Out-Null
$ws.Range("W4").PasteSpecial(-4122) | Out-Null
$ws.Range("W4").Value = "("

# X4 was ">=" -> becomes ")", same parenthesis style.
$ws.Range("A4").Copy() | Out-Null
$ws.Range("X4").PasteSpecial(-4122) | Out-Null
$ws.Range("X4").Value = ")"

# Update the saved cursor/selection position recorded in the sheet view.
$ws.Range("R10").Select() | Out-Null
